$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference style (no explicit formatting) used to restore cell style
# after temporarily forcing a text number format, so text-like numeric
# strings (e.g. "1.50", "0.0480") are preserved exactly instead of being
# auto-converted to floating point numbers by Excel.
$plainStyle = $ws.Range("B2").Style

$ws.Range("D2").Value = '25.960.06'
$ws.Range("E2").Value = '  -0.36%  '
$ws.Range("D3").Value = '1.621.78'
$ws.Range("E3").Value = '  -1.08%  '
$ws.Range("E4").Value = '  -0.42%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '212.72'
$ws.Range("D5").Style = $plainStyle
$ws.Range("E5").Value = '  -0.92%  '
$ws.Range("E6").Value = '  -1.53%  '
$ws.Range("E7").Value = '  -0.40%  '
$ws.Range("E8").Value = '  -0.95%  '
$ws.Range("E9").Value = '  -1.39%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '18.44'
$ws.Range("D10").Style = $plainStyle
$ws.Range("E10").Value = '  -1.31%  '
$ws.Range("E11").Value = '  -0.23%  '
$ws.Range("D12").Value = '1.846.20'
$ws.Range("E12").Value = '  -1.15%  '
$ws.Range("D13").Value = '1.610.82'
$ws.Range("E13").Value = '  -1.82%  '
$ws.Range("E14").Value = '  -1.60%  '
$ws.Range("E15").Value = '  -1.40%  '
$ws.Range("D16").Value = '25.969.60'
$ws.Range("E16").Value = '  -0.36%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '61.72'
$ws.Range("D17").Style = $plainStyle
$ws.Range("E17").Value = '  -1.10%  '
$ws.Range("E18").Value = '  -1.37%  '
$ws.Range("E19").Value = '  -0.40%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '191.98'
$ws.Range("D20").Style = $plainStyle
$ws.Range("E20").Value = '  +0.33%  '
$ws.Range("E21").Value = '  -0.38%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.54'
$ws.Range("D22").Style = $plainStyle
$ws.Range("E22").Value = '  -0.71%  '
$ws.Range("E23").Value = '  -1.91%  '
$ws.Range("E24").Value = '  +1.46%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '143.82'
$ws.Range("D25").Style = $plainStyle
$ws.Range("E25").Value = '  -0.19%  '
$ws.Range("E26").Value = '  -0.42%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.72'
$ws.Range("D27").Style = $plainStyle
$ws.Range("E28").Value = '  -1.99%  '
$ws.Range("E29").Value = '  -0.10%  '
$ws.Range("E30").Value = '  -1.16%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0480'
$ws.Range("D31").Style = $plainStyle
$ws.Range("E31").Value = '  -1.74%  '
$ws.Range("E32").Value = '  -1.56%  '
$ws.Range("E33").Value = '  -2.61%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.50'
$ws.Range("D34").Style = $plainStyle
$ws.Range("E34").Value = '  -0.96%  '
$ws.Range("E35").Value = '  -0.82%  '
$ws.Range("D36").Value = '1.127.09'
$ws.Range("E36").Value = '  -0.47%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.845'
$ws.Range("D37").Style = $plainStyle
$ws.Range("E37").Value = '  -3.70%  '
$ws.Range("E38").Value = '  -2.21%  '
$ws.Range("E39").Value = '  -2.09%  '
$ws.Range("E40").Value = '  -1.25%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '97.72'
$ws.Range("D41").Style = $plainStyle
$ws.Range("E41").Value = '  -1.22%  '
$ws.Range("D42").Value = '1.757.45'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.759'
$ws.Range("D43").Style = $plainStyle
$ws.Range("E43").Value = '  -3.42%  '
$ws.Range("E44").Value = '  -4.19%  '
$ws.Range("B45").Value = 'BabyDogeCoin'
$ws.Range("C45").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D45").Value = '0.0₆0113'
$ws.Range("E45").Value = '  -0.69%  '
$ws.Range("B46").Value = 'RenderToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.52'
$ws.Range("D46").Style = $plainStyle
$ws.Range("E46").Value = '  +2.03%  '
$ws.Range("B47").Value = 'Aave'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '54.17'
$ws.Range("D47").Style = $plainStyle
$ws.Range("E47").Value = '  -2.25%  '
$ws.Range("B48").Value = 'Cronos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0516'
$ws.Range("D48").Style = $plainStyle
$ws.Range("E48").Value = '  -2.05%  '
$ws.Range("B49").Value = 'Mantle'
$ws.Range("C49").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.410'
$ws.Range("D49").Style = $plainStyle
$ws.Range("E49").Value = '  -1.12%  '
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.50'
$ws.Range("D50").Style = $plainStyle
$ws.Range("E50").Value = '  -0.80%  '
$ws.Range("B51").Value = 'USDD'
$ws.Range("C51").Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.00'
$ws.Range("D51").Style = $plainStyle
$ws.Range("E51").Value = '  -0.23%  '
